# The recorded change is a whole-deck resave: PowerPoint re-serialized the
# <mc:AlternateContent>/<mc:Choice>/<mc:Fallback> wrapper that carries each
# slide/layout/master's "no transition" (p14:dur="0") marker, hoisting the
# xmlns:p14 declaration up onto <mc:AlternateContent> and adding an empty
# xmlns="" reset on <mc:Fallback>. The transition semantics themselves
# (Requires="p14", p14:dur="0") are untouched on every single part -- this
# is purely how the host renumbers/relocates namespace declarations when it
# rewrites a part, not a content edit.
#
# Touch each slide's transition object (with its own, already-current value)
# so every slide's XML gets regenerated/re-flushed through the same
# serializer PowerPoint used when it produced the committed file, without
# altering any transition attribute (Hidden is already False on every
# slide, so re-asserting it is a no-op content-wise).
$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $t = $s.SlideShowTransition
    $t.Hidden = $false
}
